# Update "Login Info" sheet with a new "Jenkins:" column (E) and select E3,
# per commit "Update excel with Jenkins".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Info")

# Add new header label in E1 and its value (1) in E2, extending the table
# that already has URL:/Username:/Password:/Adjustresolution: in A1:D1.
$ws.Range("E1").Value = "Jenkins:"
$ws.Range("E2").Value = 1

# Update the active selection on the sheet to E3 (matches saved view state).
$ws.Range("E3").Select() | Out-Null
